# Adds two new annotation columns (K: KLIFS residue id, L: free-text
# comment) to the "Counts" sheet, populated for a subset of rows, plus a
# header for each ("5dls" / "comment").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Counts")

# Seed K1/L1 with the same header style as the existing header cells
# (bold, rotated, centered) by copying J1's formatting, then overwrite
# the copied value with the real header text.
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Range("J1").Copy($ws.Range("L1"))
$ws.Range("K1").Value = '5dls'
$ws.Range("L1").Value = 'comment'

$ws.Range("K2").Value = 'Q'

$ws.Range("K3").Value = 'T'

$ws.Range("K4").Value = 'L15'
$ws.Range("L4").Value = 'ok'

$ws.Range("K5").Value = 'G'

$ws.Range("K6").Value = 'E'

$ws.Range("K7").Value = 'G'

$ws.Range("K8").Value = 'A36'

$ws.Range("K9").Value = 'Y'
$ws.Range("L9").Value = 'very dependent on the conformation'

$ws.Range("K10").Value = 'G'

$ws.Range("K11").Value = 'E'

$ws.Range("K12").Value = 'V23'
$ws.Range("L12").Value = 'ok'

$ws.Range("K13").Value = 'Q'

$ws.Range("K14").Value = 'L15'

$ws.Range("K15").Value = 'V23'

$ws.Range("K16").Value = 'A36'
$ws.Range("L16").Value = 'ok'

$ws.Range("K17").Value = 'V23'

$ws.Range("K18").Value = 'K38'
$ws.Range("L18").Value = 'ok'

$ws.Range("K19").Value = 'I'

$ws.Range("K20").Value = 'V23'

$ws.Range("K21").Value = 'N'

$ws.Range("K24").Value = 'K'

$ws.Range("K25").Value = 'E55'
$ws.Range("L25").Value = 'a bit far away for type I no?'

$ws.Range("K26").Value = 'I'

$ws.Range("K28").Value = 'I'

$ws.Range("K29").Value = 'N59'
$ws.Range("L29").Value = 'a bit far away for type I no?'

$ws.Range("K30").Value = 'K'

$ws.Range("K36").Value = 'V'

$ws.Range("K37").Value = 'V68'
$ws.Range("L37").Value = 'ok'

$ws.Range("K38").Value = 'K'

$ws.Range("K45").Value = 'F83'

$ws.Range("K46").Value = 'L84'
$ws.Range("L46").Value = 'ok'

$ws.Range("K47").Value = 'E85'
$ws.Range("L47").Value = 'I don''t understand the sidechain interactions -> examples of structures here?'

$ws.Range("K48").Value = 'Y86'
$ws.Range("L48").Value = 'same, rather not directly in contact with sidechain'

$ws.Range("K49").Value = 'C87'
$ws.Range("L49").Value = 'same, rather not directly in contact with sidechain'

$ws.Range("K50").Value = 'S88'
$ws.Range("L50").Value = 'same, example where the sidechain is in contact?'

# Match the author's final cursor position/selection in the saved file.
$ws.Range("K52").Select()
